$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manager credentials shown in row 2
$ws.Range("A2").Value = "mngr302125"
$ws.Range("B2").Value = "pAjapEq"

# Move the selection to B2 (matches the saved selection in the sheet view)
$ws.Range("B2").Select()
